$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 96 (pushes the existing rows 96-105 down to 97-106)
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with a new price record
$ws.Cells.Item(96, 1).Value = 7
$ws.Cells.Item(96, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(96, 3).Value = "Ñuble"
$ws.Cells.Item(96, 4).Value = 45127
$ws.Cells.Item(96, 5).Value = 16
$ws.Cells.Item(96, 6).Value = 100112013
$ws.Cells.Item(96, 7).Value = "Alcachofa"
$ws.Cells.Item(96, 8).Value = "Argentina(o)"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 30
$ws.Cells.Item(96, 11).Value = 17000
$ws.Cells.Item(96, 12).Value = 17000
$ws.Cells.Item(96, 13).Value = 17000
$ws.Cells.Item(96, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(96, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(96, 16).Value = 340
$ws.Cells.Item(96, 17).Value = 50
$ws.Cells.Item(96, 18).Value = "Hortaliza"
